$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Joshua Aguirre's availability (row 4)
$ws.Range("B4").Value = "5pm-MN"
$ws.Range("C4").Value = "5pm-MN"
$ws.Range("D4").Value = "5pm-MN"
$ws.Range("E4").Value = "5pm-MN"
$ws.Range("F4").Value = "5pm-MN"
$ws.Range("G4").Value = "8am-MN"
$ws.Range("H4").Value = "8am-MN"

# Zachary Hickerson's availability (row 5)
$ws.Range("B5").Value = "12pm-5pm"
$ws.Range("C5").Value = "2pm-5pm"
$ws.Range("D5").Value = "12pm-5pm"
$ws.Range("E5").Value = "2pm-4pm"
$ws.Range("F5").Value = "12pm-MN"
$ws.Range("G5").Value = "6pm-10pm"
$ws.Range("H5").Value = "6pm-MN"

# Clear Adam Furbee's old placeholder availability (row 10)
$ws.Range("B10:H10").ClearContents()

# Update active selection to D6
$ws.Range("D6").Select()
